# Update 28 Mai 2025
# Applies the data refresh captured in the commit "Update 28 Mai 2025" to the
# genre / commune / trimestre statistics workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Genre_Individuel"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Genre_Individuel")
$ws.Range("B2").Value = 10528
$ws.Range("C2").Value = 83.146422366134885
$ws.Range("B3").Value = 2134
$ws.Range("C3").Value = 16.853577633865111

# ---------------------------------------------------------------------------
# Sheet "Genre_Collectif"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Genre_Collectif")
$ws.Range("B2").Value = 28042
$ws.Range("C2").Value = 81.41330855882012
$ws.Range("B3").Value = 6402
$ws.Range("C3").Value = 18.586691441179891

# ---------------------------------------------------------------------------
# Sheet "Genre_Total"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Genre_Total")
$ws.Range("B2").Value = 38570
$ws.Range("C2").Value = 81.879166135948708
$ws.Range("B3").Value = 8536
$ws.Range("C3").Value = 18.120833864051288

# ---------------------------------------------------------------------------
# Sheet "Genre_Mandataires"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Genre_Mandataires")
$ws.Range("B2").Value = 5798
$ws.Range("C2").Value = 94.661224489795913
$ws.Range("B3").Value = 327
$ws.Range("C3").Value = 5.3387755102040817

# ---------------------------------------------------------------------------
# Sheet "Récapitulatif_Genre" (mirrors the four sheets above)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Récapitulatif_Genre")
$ws.Range("B2").Value = 10528
$ws.Range("C2").Value = 83.146422366134885
$ws.Range("D2").Value = 28042
$ws.Range("E2").Value = 81.41330855882012
$ws.Range("F2").Value = 38570
$ws.Range("G2").Value = 81.879166135948708
$ws.Range("H2").Value = 5798
$ws.Range("I2").Value = 94.661224489795913

$ws.Range("B3").Value = 2134
$ws.Range("C3").Value = 16.853577633865111
$ws.Range("D3").Value = 6402
$ws.Range("E3").Value = 18.586691441179891
$ws.Range("F3").Value = 8536
$ws.Range("G3").Value = 18.120833864051288
$ws.Range("H3").Value = 327
$ws.Range("I3").Value = 5.3387755102040817

$ws.Range("A1:I3").Select()

# ---------------------------------------------------------------------------
# Sheet "Analyse_Commune" -- dataset swapped from commune names to commune
# codes, and two extra rows appended.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Analyse_Commune")

# Grow the table by copying the formatting of the last two existing rows
# down into rows 11:12 before rewriting every row's content.
$ws.Range("A9:F10").Copy($ws.Range("A11:F12"))

$communeRows = @(
    @(2,  5120301, 450,  5060, 5510, 8.1669691470054442, 91.833030852994554),
    @(3,  5120302, 200,  3815, 4015, 4.9813200498132,    95.018679950186808),
    @(4,  5120303, 49,   953,  1002, 4.8902195608782426, 95.109780439121764),
    @(5,  5220202, 97,   1081, 1178, 8.2342954159592523, 91.765704584040748),
    @(6,  5220302, 1015, 5076, 6091, 16.66393038909867,  83.336069610901333),
    @(7,  5220303, 484,  6427, 6911, 7.0033280277817962, 92.996671972218209),
    @(8,  13120101, 2381, 5784, 8165, 29.161053276178809, 70.838946723821181),
    @(9,  13120103, 363,  1167, 1530, 23.725490196078429, 76.274509803921561),
    @(10, 13120201, 1353, 3504, 4857, 27.856701667696111, 72.1432983323039),
    @(11, 13120202, 1344, 3041, 4385, 30.649942987457241, 69.350057012542749),
    @(12, 13320102, 800,  2650, 3450, 23.188405797101449, 76.811594202898547)
)

foreach ($row in $communeRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
}

$ws.Range("A1:F12").Select()

# ---------------------------------------------------------------------------
# Sheet "Analyse_Trimestre"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Analyse_Trimestre")
$ws.Range("C3").Value = 18676
$ws.Range("D3").Value = 22715
$ws.Range("E3").Value = 17.78120184899846
$ws.Range("F3").Value = 82.218798151001536

$ws.Range("B4").Value = 2890
$ws.Range("C4").Value = 13830
$ws.Range("D4").Value = 16720
$ws.Range("E4").Value = 17.284688995215308
$ws.Range("F4").Value = 82.715311004784681

$ws.Range("A1:F4").Select()
